$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44685
$ws.Range("J3").Value = 400

# Row 4
$ws.Range("D4").Value = 44827
$ws.Range("J4").Value = 1200
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2250
$ws.Range("P4").Value = 750

# Row 5
$ws.Range("D5").Value = 44883
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 1800
$ws.Range("M5").Value = 1900
$ws.Range("P5").Value = 633

# Row 6
$ws.Range("D6").Value = 44848
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 1500
$ws.Range("M6").Value = 1750
$ws.Range("P6").Value = 583

# Row 7
$ws.Range("D7").Value = 44881
$ws.Range("K7").Value = 1900
$ws.Range("M7").Value = 1950
$ws.Range("P7").Value = 650

# Row 8
$ws.Range("D8").Value = 44911
$ws.Range("J8").Value = 700
$ws.Range("K8").Value = 1800
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 1900
$ws.Range("P8").Value = 633
